$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header in F1, matching style of existing headers (e.g. E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("F1").Style = $ws.Range("E1").Style

# Fill time_taken values for each data row
$ws.Range("F2").Value = "2021-10-05 13:41:34.843038"
$ws.Range("F3").Value = "2021-10-05 13:41:34.843050"
$ws.Range("F4").Value = "2021-10-05 13:41:34.843053"
$ws.Range("F5").Value = "2021-10-05 13:41:34.843057"
$ws.Range("F6").Value = "2021-10-05 13:41:34.843060"
$ws.Range("F7").Value = "2021-10-05 13:41:34.843063"
$ws.Range("F8").Value = "2021-10-05 13:41:34.843066"
